$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - value + the same look (bold/centered/bordered) as the
# existing header row (B1:G1), obtained by copying G1's formatting over.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New "Save" column values for the existing data rows (2-5)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
